# Update cryptos list (prices in column D, 1h volume/change in column E)
# Commit: "Updated cryptos list on Mon Feb  5 17:52:21 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new Price (column D) / new Volume(1h) (column E) values.
# $null for D means the price text is unchanged for that row.
$updates = @(
    @{ Row = 2;  D = "42.749.18"; E = "  -0.42%  " },
    @{ Row = 3;  D = "2.303.53";  E = "  -0.16%  " },
    @{ Row = 4;  D = "1.00";      E = "  -0.02%  " },
    @{ Row = 5;  D = "301.17";    E = "  -1.52%  " },
    @{ Row = 6;  D = "96.13";     E = "  -1.26%  " },
    @{ Row = 7;  D = "0.510";     E = "  -0.31%  " },
    @{ Row = 8;  D = $null;       E = "  +0.07%  " },
    @{ Row = 9;  D = $null;       E = "  -1.82%  " },
    @{ Row = 10; D = "34.77";     E = "  -2.64%  " },
    @{ Row = 11; D = "19.30";     E = "  +5.04%  " },
    @{ Row = 12; D = $null;       E = "  -0.94%  " },
    @{ Row = 13; D = $null;       E = "  -0.35%  " },
    @{ Row = 14; D = "6.79";      E = "  +0.20%  " },
    @{ Row = 15; D = "2.653.34";  E = "  -0.39%  " },
    @{ Row = 16; D = "2.293.62";  E = "  -0.26%  " },
    @{ Row = 17; D = $null;       E = "  +0.24%  " },
    @{ Row = 18; D = "42.692.89"; E = "  -0.34%  " },
    @{ Row = 19; D = "12.35";     E = "  -5.99%  " },
    @{ Row = 20; D = $null;       E = "  -1.39%  " },
    @{ Row = 21; D = "6.02";      E = "  -0.52%  " },
    @{ Row = 22; D = "67.93";     E = "  +0.38%  " },
    @{ Row = 23; D = $null;       E = "  +3.94%  " },
    @{ Row = 24; D = "235.22";    E = "  -0.78%  " },
    @{ Row = 25; D = $null;       E = "  +0.17%  " },
    @{ Row = 26; D = $null;       E = "  -3.18%  " },
    @{ Row = 27; D = "24.61";     E = "  -3.88%  " },
    @{ Row = 28; D = $null;       E = "  -0.61%  " },
    @{ Row = 29; D = "164.83";    E = "  -1.59%  " },
    @{ Row = 30; D = "9.07";      E = "  -0.28%  " },
    @{ Row = 31; D = "32.19";     E = "  -2.75%  " },
    @{ Row = 32; D = $null;       E = "  -0.02%  " },
    @{ Row = 33; D = "4.97";      E = "  -0.84%  " },
    @{ Row = 34; D = "17.43";     E = "  -0.25%  " },
    @{ Row = 35; D = "4.43";      E = "  -7.45%  " },
    @{ Row = 36; D = $null;       E = "  +0.91%  " },
    @{ Row = 37; D = $null;       E = "  -2.89%  " },
    @{ Row = 38; D = $null;       E = "  -2.01%  " },
    @{ Row = 39; D = $null;       E = "  -0.57%  " },
    @{ Row = 40; D = "2.71";      E = "  -0.69%  " },
    @{ Row = 41; D = "0.108";     E = "  -1.49%  " },
    @{ Row = 42; D = "19.76";     E = "  +5.92%  " },
    @{ Row = 43; D = "1.972.06";  E = "  -1.77%  " },
    @{ Row = 44; D = $null;       E = "  +4.43%  " },
    @{ Row = 45; D = $null;       E = "  -0.75%  " },
    @{ Row = 46; D = $null;       E = "  -3.27%  " },
    @{ Row = 47; D = "2.76";      E = "  -0.78%  " },
    @{ Row = 48; D = "2.85";      E = "  -3.31%  " },
    @{ Row = 49; D = "2.527.15";  E = "  -0.20%  " },
    @{ Row = 50; D = "53.16";     E = "  -1.61%  " },
    @{ Row = 51; D = "71.55";     E = "  -0.71%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # The price column holds plain text (e.g. "42.749.18", "1.00") that looks
        # numeric. Assigning a Value directly would make Excel reinterpret it as a
        # number (and drop things like trailing zeros or the extra "thousands" dot).
        # Writing it as a quoted text formula and then collapsing the formula down
        # to its literal value via Copy/PasteSpecial(values) keeps it a genuine
        # text cell without touching any cell formatting/styles.
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.Formula = '="' + $u.D + '"'
        $cell.Copy() | Out-Null
        $cell.PasteSpecial(-4163) | Out-Null
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

$excel.CutCopyMode = 0
